$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.890.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.041.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.81'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.30%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.91'
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.85%  '
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.340.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.807'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.040.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.866.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +16.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").Value = '  +4.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0619'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0870'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.110'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.49%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.97%  '
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +14.28%  '
$ws.Range("E47").Value = '  +5.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.282.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.226.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.31%  '
